$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28..133 down to 29..134
$ws.Rows("28:28").Insert()

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44565
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 30000
$ws.Range("N28").Value = "`$/saco 25 kilos"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1200
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
